{"js": "// Word (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Change 1: remove the two leading empty \"Titre2\"/centered paragraphs at the\n//           very top of the document (the third \u2014 with the actual title text\n//           \"Authorization Form for Access to Controlled Area\" \u2014 is kept).\n// Change 2: rework the \"Signature of employer\" block:\n//             - \"Signature of employer: \" (+ trailing run \" \") -> \"Signature of employer:\"\n//               (single run, no trailing space, plain justified paragraph)\n//             - new paragraph holding the literal placeholder \"[SignatureField#1]\"\n//             - the following empty spacer paragraph becomes a plain centered\n//               empty paragraph (its old rFonts/sz formatting is dropped)\n//             - the page-break paragraph's pPr is cleared\n//             - the following empty \"Normal\" spacer paragraph is dropped\n//               (merged away), and the final \"AGLAE Controlled Area\" heading\n//               paragraph becomes bold (w:b/w:bCs) and centered, with an\n//               explicit (empty) run-level rPr.\n\nconst body = context.document.body;\n\n// ---- Change 1: drop the two leading empty Titre2/centered paragraphs ----\nconst leadParas = body.paragraphs;\nleadParas.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (const p of leadParas.items) {\n  if (p.style === \"Heading 2\" && p.text === \"\") {\n    p.delete();\n  } else {\n    break;\n  }\n}\nawait context.sync();\n\n// ---- Change 2: locate the \"Signature of employer\" paragraph and rebuild\n//      the 5-paragraph block that follows (and includes) it. ----\nconst results = body.search(\"Signature of employer:\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst startPara = results.items[0].paragraphs.getFirst();\nstartPara.load(\"items\");\nconst allParas = body.paragraphs;\nallParas.load(\"items\");\nawait context.sync();\n\n// Find the index of the paragraph containing \"Signature of employer:\" so we\n// can grab the 4 paragraphs that immediately follow it too.\nlet idx = -1;\nfor (let i = 0; i < allParas.items.length; i++) {\n  if (allParas.items[i].text.indexOf(\"Signature of employer:\") === 0) {\n    idx = i;\n    break;\n  }\n}\n\nconst firstPara = allParas.items[idx];\nconst lastPara = allParas.items[idx + 4]; // the \"AGLAE Controlled Area\" paragraph\n\nconst blockRange = firstPara.getRange(\"Start\").expandTo(lastPara.getRange(\"End\"));\n\nconst newBodyXml =\n  '<w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:rPr/><w:t>Signature of employer:</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:rPr/><w:t>[SignatureField#1]</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' +\n  '<w:p><w:pPr/><w:r><w:br w:type=\"page\"/></w:r></w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:cs=\"Arial\"/><w:b w:val=\"1\"/><w:bCs w:val=\"1\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>AGLAE Controlled Area</w:t></w:r></w:p>';\n\nconst packageXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newBodyXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nblockRange.insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop script \u2014 $word / $d (ActiveDocument) are pre-seeded.\n#\n# Change 1: remove the two leading empty \"Heading 2\" (Titre2)/centered\n#           paragraphs at the very top of the document (the third one,\n#           holding the actual title \"Authorization Form for Access to\n#           Controlled Area\", is left untouched).\n# Change 2: rework the \"Signature of employer\" block (5 paragraphs):\n#             - \"Signature of employer: \" (+ trailing run \" \") -> \"Signature of employer:\"\n#             - new paragraph holding the literal placeholder \"[SignatureField#1]\"\n#             - the following empty spacer paragraph becomes a plain centered\n#               empty paragraph (old rFonts/sz formatting dropped)\n#             - the page-break paragraph's pPr is cleared\n#             - the following empty \"Normal\" spacer paragraph is dropped\n#               (merged away), and the final \"AGLAE Controlled Area\" heading\n#               paragraph becomes bold + centered with an explicit empty\n#               run-level rPr.\n\n$d = $word.ActiveDocument\n\n# ---- Change 1: drop the two leading empty Heading2/centered paragraphs ----\nwhile ($true) {\n    $p = $d.Paragraphs.Item(1)\n    $plainText = $p.Range.Text -replace \"`r$\", \"\"\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $plainText -eq \"\") {\n        $p.Range.Delete()\n    } else {\n        break\n    }\n}\n\n# ---- Change 2: locate \"Signature of employer:\" and rebuild the block ----\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"Signature of employer:\")\n$firstPara = $findRange.Paragraphs(1)\n\n$allParas = $d.Paragraphs\n$startIdx = -1\nfor ($i = 1; $i -le $allParas.Count; $i++) {\n    if ($allParas.Item($i).Range.Start -eq $firstPara.Range.Start) {\n        $startIdx = $i\n        break\n    }\n}\n\n$lastPara = $allParas.Item($startIdx + 4)   # the \"AGLAE Controlled Area\" paragraph\n$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n\n$newBodyXml = '<w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:rPr/><w:t>Signature of employer:</w:t></w:r></w:p>' + `\n    '<w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:rPr/><w:t>[SignatureField#1]</w:t></w:r></w:p>' + `\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' + `\n    '<w:p><w:pPr/><w:r><w:br w:type=\"page\"/></w:r></w:p>' + `\n    '<w:p><w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:cs=\"Arial\"/><w:b w:val=\"1\"/><w:bCs w:val=\"1\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>AGLAE Controlled Area</w:t></w:r></w:p>'\n\n$packageXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + $newBodyXml + '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n    '</pkg:package>'\n\n$blockRange.InsertXML($packageXml)\n"}
